# Adds Rating / Total Ratings / Popularity columns (O,P,Q) to the Report
# sheet, backfills them for every existing product row, and inserts two
# new product rows ("Mix Veg  Soup" and "chicken Sour") in the middle of
# the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New columns O/P/Q (widths match the existing N column formatting)
# ---------------------------------------------------------------------
$ws.Columns.Item(15).ColumnWidth = 10.7109375
$ws.Columns.Item(16).ColumnWidth = 10.7109375
$ws.Columns.Item(17).ColumnWidth = 15.7109375

$ws.Range("O1").Value = "Rating"
$ws.Range("P1").Value = "Total Ratings"
$ws.Range("Q1").Value = "Popularity"

# ---------------------------------------------------------------------
# 2. Rating data for the 7 original product rows (rows 2-8)
# ---------------------------------------------------------------------
$ratings = @{
    2 = @("4",   "0", "2")
    3 = @("4",   "0", "0")
    4 = @("4",   "0", "2")
    5 = @("4.5", "0", "0")
    6 = @("4.5", "0", "3")
    7 = @("4",   "0", "0")
    8 = @("4.5", "0", "0")
}
foreach ($r in $ratings.Keys) {
    $vals = $ratings[$r]
    $ws.Range("O$r").Value = $vals[0]
    $ws.Range("P$r").Value = $vals[1]
    $ws.Range("Q$r").Value = $vals[2]
}

# ---------------------------------------------------------------------
# 3. Insert the "Mix Veg  Soup" row as the new row 9 (old rows 9-13 slide
#    down to 10-14).
# ---------------------------------------------------------------------
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "veg"
$ws.Range("B9").Value = "e70158bb-c576-4135-89f2-f007a943a58f"
$ws.Range("C9").Value = "Mix Veg  Soup"
$ws.Range("D9").Value = "Veg"
$ws.Range("E9").Value = "15 min"
$ws.Range("F9").Value = "0"
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = ""
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = "Soup is a primarily liquid food, generally served warm or hot, that is made by combining ingredients of meat or vegetables with stock, or water. Hot soups are additionally characterized by boiling solid ingredients in liquids in a pot until the flavors are extracted, forming a broth"
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = "cutlery"
$ws.Range("N9").Value = "Active"
$ws.Range("O9").Value = "4"
$ws.Range("P9").Value = "0"
$ws.Range("Q9").Value = "0"

# ---------------------------------------------------------------------
# 4. Ratings for the rows that used to be 9-12 and are now 10-12 (the
#    "Paneer Samosa" / old row 12 and "Sweet Corn Soup" / old row 13 get
#    handled after the second insert below).
# ---------------------------------------------------------------------
$ratings2 = @{
    10 = @("4", "0", "1")
    11 = @("4", "0", "0")
    12 = @("4", "0", "0")
}
foreach ($r in $ratings2.Keys) {
    $vals = $ratings2[$r]
    $ws.Range("O$r").Value = $vals[0]
    $ws.Range("P$r").Value = $vals[1]
    $ws.Range("Q$r").Value = $vals[2]
}

# ---------------------------------------------------------------------
# 5. Insert the "chicken Sour" row as the new row 13 (old rows 13-14,
#    i.e. "Paneer Samosa" and "Sweet Corn Soup", slide down to 14-15).
# ---------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = "Soups"
$ws.Range("B13").Value = "aa7b0190-97b1-4519-bf66-df04f5f33b85"
$ws.Range("C13").Value = "chicken Sour"
$ws.Range("D13").Value = "Non-Veg"
$ws.Range("E13").Value = "20 min"
$ws.Range("F13").Value = "260"
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = ""
$ws.Range("J13").Value = 260
$ws.Range("K13").Value = "Soup is a primarily liquid food, generally served warm or hot, that is made by combining ingredients of meat or vegetables with stock, or water. Hot soups are additionally characterized by boiling solid ingredients in liquids in a pot until the flavors are extracted, forming a broth"
$ws.Range("L13").Value = ""
$ws.Range("N13").Value = "Active"
$ws.Range("O13").Value = "0"
$ws.Range("P13").Value = "0"
$ws.Range("Q13").Value = "0"

# ---------------------------------------------------------------------
# 6. Ratings for the two trailing rows (old rows 12 & 13, "Paneer Samosa"
#    and "Sweet Corn Soup", now sitting at 14 & 15).
# ---------------------------------------------------------------------
$ws.Range("O14").Value = "4"
$ws.Range("P14").Value = "0"
$ws.Range("Q14").Value = "0"

$ws.Range("O15").Value = "4"
$ws.Range("P15").Value = "0"
$ws.Range("Q15").Value = "0"
